$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Charvarius Ward'
$ws.Cells.Item(2, 2).Value = 'Group1'
$ws.Cells.Item(2, 3).Value = 8.666666666666666
$ws.Cells.Item(2, 4).Value = 64
$ws.Cells.Item(2, 5).Value = 47.33333333333334
$ws.Cells.Item(2, 6).Value = 16.66666666666667

$ws.Cells.Item(3, 1).Value = 'Charvarius Ward'
$ws.Cells.Item(3, 2).Value = 'Group2'
$ws.Cells.Item(3, 3).Value = 13.66666666666667
$ws.Cells.Item(3, 4).Value = 71
$ws.Cells.Item(3, 5).Value = 51
$ws.Cells.Item(3, 6).Value = 20

$ws.Cells.Item(4, 1).Value = 'Charvarius Ward'
$ws.Cells.Item(4, 2).Value = 'Difference'
$ws.Cells.Item(4, 3).Value = 5
$ws.Cells.Item(4, 4).Value = 7
$ws.Cells.Item(4, 5).Value = 3.666666666666664
$ws.Cells.Item(4, 6).Value = 3.333333333333332

$ws.Cells.Item(5, 1).Value = 'D.J. Reed'
$ws.Cells.Item(5, 2).Value = 'Group1'
$ws.Cells.Item(5, 3).Value = 6.333333333333333
$ws.Cells.Item(5, 4).Value = 51
$ws.Cells.Item(5, 5).Value = 39.66666666666666
$ws.Cells.Item(5, 6).Value = 11.33333333333333

$ws.Cells.Item(6, 1).Value = 'D.J. Reed'
$ws.Cells.Item(6, 2).Value = 'Group2'
$ws.Cells.Item(6, 3).Value = 10.66666666666667
$ws.Cells.Item(6, 4).Value = 73.33333333333333
$ws.Cells.Item(6, 5).Value = 59
$ws.Cells.Item(6, 6).Value = 14.33333333333333

$ws.Cells.Item(7, 1).Value = 'D.J. Reed'
$ws.Cells.Item(7, 2).Value = 'Difference'
$ws.Cells.Item(7, 3).Value = 4.333333333333333
$ws.Cells.Item(7, 4).Value = 22.33333333333333
$ws.Cells.Item(7, 5).Value = 19.33333333333334
$ws.Cells.Item(7, 6).Value = 3

$ws.Cells.Item(8, 1).Value = 'Denzel Ward'
$ws.Cells.Item(8, 2).Value = 'Group1'
$ws.Cells.Item(8, 3).Value = 13
$ws.Cells.Item(8, 4).Value = 44.33333333333334
$ws.Cells.Item(8, 5).Value = 36.66666666666666
$ws.Cells.Item(8, 6).Value = 7.666666666666667

$ws.Cells.Item(9, 1).Value = 'Denzel Ward'
$ws.Cells.Item(9, 2).Value = 'Group2'
$ws.Cells.Item(9, 3).Value = 15
$ws.Cells.Item(9, 4).Value = 45.33333333333334
$ws.Cells.Item(9, 5).Value = 37.66666666666666
$ws.Cells.Item(9, 6).Value = 7.666666666666667

$ws.Cells.Item(10, 1).Value = 'Denzel Ward'
$ws.Cells.Item(10, 2).Value = 'Difference'
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0

$ws.Cells.Item(11, 1).Value = 'Isaac Yiadom'
$ws.Cells.Item(11, 2).Value = 'Group1'
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 33
$ws.Cells.Item(11, 5).Value = 23.33333333333333
$ws.Cells.Item(11, 6).Value = 9.666666666666666

$ws.Cells.Item(12, 1).Value = 'Isaac Yiadom'
$ws.Cells.Item(12, 2).Value = 'Group2'
$ws.Cells.Item(12, 3).Value = 6.666666666666667
$ws.Cells.Item(12, 4).Value = 28.11111111111111
$ws.Cells.Item(12, 5).Value = 21.22222222222222
$ws.Cells.Item(12, 6).Value = 6.888888888888888

$ws.Cells.Item(13, 1).Value = 'Isaac Yiadom'
$ws.Cells.Item(13, 2).Value = 'Difference'
$ws.Cells.Item(13, 3).Value = 3.666666666666667
$ws.Cells.Item(13, 4).Value = -4.888888888888886
$ws.Cells.Item(13, 5).Value = -2.111111111111107
$ws.Cells.Item(13, 6).Value = -2.777777777777778

$ws.Cells.Item(14, 1).Value = 'J.T. Gray'
$ws.Cells.Item(14, 2).Value = 'Group1'
$ws.Cells.Item(14, 3).Value = 0.3333333333333333
$ws.Cells.Item(14, 4).Value = 15.66666666666667
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 2.666666666666667

$ws.Cells.Item(15, 1).Value = 'J.T. Gray'
$ws.Cells.Item(15, 2).Value = 'Group2'
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 18
$ws.Cells.Item(15, 5).Value = 9.333333333333334
$ws.Cells.Item(15, 6).Value = 8.666666666666666

$ws.Cells.Item(16, 1).Value = 'J.T. Gray'
$ws.Cells.Item(16, 2).Value = 'Difference'
$ws.Cells.Item(16, 3).Value = -0.3333333333333333
$ws.Cells.Item(16, 4).Value = 2.333333333333334
$ws.Cells.Item(16, 5).Value = -3.666666666666666
$ws.Cells.Item(16, 6).Value = 6

$ws.Cells.Item(17, 1).Value = 'Kendall Fuller'
$ws.Cells.Item(17, 2).Value = 'Group1'
$ws.Cells.Item(17, 3).Value = 9.666666666666666
$ws.Cells.Item(17, 4).Value = 58.66666666666666
$ws.Cells.Item(17, 5).Value = 45.33333333333334
$ws.Cells.Item(17, 6).Value = 13.33333333333333

$ws.Cells.Item(18, 1).Value = 'Kendall Fuller'
$ws.Cells.Item(18, 2).Value = 'Group2'
$ws.Cells.Item(18, 3).Value = 9.666666666666666
$ws.Cells.Item(18, 4).Value = 60.33333333333334
$ws.Cells.Item(18, 5).Value = 42.66666666666666
$ws.Cells.Item(18, 6).Value = 17.66666666666667

$ws.Cells.Item(19, 1).Value = 'Kendall Fuller'
$ws.Cells.Item(19, 2).Value = 'Difference'
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 1.666666666666671
$ws.Cells.Item(19, 5).Value = -2.666666666666671
$ws.Cells.Item(19, 6).Value = 4.333333333333334

$ws.Cells.Item(20, 1).Value = 'Kevin Byard'
$ws.Cells.Item(20, 2).Value = 'Group1'
$ws.Cells.Item(20, 3).Value = 9.666666666666666
$ws.Cells.Item(20, 4).Value = 94.33333333333333
$ws.Cells.Item(20, 5).Value = 65.33333333333333
$ws.Cells.Item(20, 6).Value = 29

$ws.Cells.Item(21, 1).Value = 'Kevin Byard'
$ws.Cells.Item(21, 2).Value = 'Group2'
$ws.Cells.Item(21, 3).Value = 5
$ws.Cells.Item(21, 4).Value = 106.4444444444444
$ws.Cells.Item(21, 5).Value = 67.44444444444444
$ws.Cells.Item(21, 6).Value = 39

$ws.Cells.Item(22, 1).Value = 'Kevin Byard'
$ws.Cells.Item(22, 2).Value = 'Difference'
$ws.Cells.Item(22, 3).Value = -4.666666666666666
$ws.Cells.Item(22, 4).Value = 12.11111111111111
$ws.Cells.Item(22, 5).Value = 2.111111111111114
$ws.Cells.Item(22, 6).Value = 10

$ws.Cells.Item(23, 1).Value = 'Michael Davis'
$ws.Cells.Item(23, 2).Value = 'Group1'
$ws.Cells.Item(23, 3).Value = 11.66666666666667
$ws.Cells.Item(23, 4).Value = 52.33333333333334
$ws.Cells.Item(23, 5).Value = 40.33333333333334
$ws.Cells.Item(23, 6).Value = 12

$ws.Cells.Item(24, 1).Value = 'Michael Davis'
$ws.Cells.Item(24, 2).Value = 'Group2'
$ws.Cells.Item(24, 3).Value = 8.666666666666666
$ws.Cells.Item(24, 4).Value = 46
$ws.Cells.Item(24, 5).Value = 35.66666666666666
$ws.Cells.Item(24, 6).Value = 10.33333333333333

$ws.Cells.Item(25, 1).Value = 'Michael Davis'
$ws.Cells.Item(25, 2).Value = 'Difference'
$ws.Cells.Item(25, 3).Value = -3
$ws.Cells.Item(25, 4).Value = -6.333333333333336
$ws.Cells.Item(25, 5).Value = -4.666666666666671
$ws.Cells.Item(25, 6).Value = -1.666666666666666

$ws.Cells.Item(26, 1).Value = 'Mike Ford'
$ws.Cells.Item(26, 2).Value = 'Group1'
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 4).Value = 10
$ws.Cells.Item(26, 5).Value = 8.666666666666666
$ws.Cells.Item(26, 6).Value = 1.333333333333333

$ws.Cells.Item(27, 1).Value = 'Mike Ford'
$ws.Cells.Item(27, 2).Value = 'Group2'
$ws.Cells.Item(27, 3).Value = 0.6666666666666666
$ws.Cells.Item(27, 4).Value = 18.66666666666667
$ws.Cells.Item(27, 5).Value = 13.33333333333333
$ws.Cells.Item(27, 6).Value = 5.333333333333333

$ws.Cells.Item(28, 1).Value = 'Mike Ford'
$ws.Cells.Item(28, 2).Value = 'Difference'
$ws.Cells.Item(28, 3).Value = -0.3333333333333334
$ws.Cells.Item(28, 4).Value = 8.666666666666668
$ws.Cells.Item(28, 5).Value = 4.666666666666668
$ws.Cells.Item(28, 6).Value = 4

$ws.Cells.Item(29, 1).Value = 'Rock Ya-Sin'
$ws.Cells.Item(29, 2).Value = 'Group1'
$ws.Cells.Item(29, 3).Value = 6.666666666666667
$ws.Cells.Item(29, 4).Value = 46
$ws.Cells.Item(29, 5).Value = 39
$ws.Cells.Item(29, 6).Value = 7

$ws.Cells.Item(30, 1).Value = 'Rock Ya-Sin'
$ws.Cells.Item(30, 2).Value = 'Group2'
$ws.Cells.Item(30, 3).Value = 3.666666666666667
$ws.Cells.Item(30, 4).Value = 20.33333333333333
$ws.Cells.Item(30, 5).Value = 16.33333333333333
$ws.Cells.Item(30, 6).Value = 4

$ws.Cells.Item(31, 1).Value = 'Rock Ya-Sin'
$ws.Cells.Item(31, 2).Value = 'Difference'
$ws.Cells.Item(31, 3).Value = -3
$ws.Cells.Item(31, 4).Value = -25.66666666666667
$ws.Cells.Item(31, 5).Value = -22.66666666666667
$ws.Cells.Item(31, 6).Value = -3

$ws.Cells.Item(32, 1).Value = 'Stephon Gilmore'
$ws.Cells.Item(32, 2).Value = 'Group1'
$ws.Cells.Item(32, 3).Value = 8.333333333333334
$ws.Cells.Item(32, 4).Value = 35.33333333333334
$ws.Cells.Item(32, 5).Value = 29.66666666666667
$ws.Cells.Item(32, 6).Value = 5.666666666666667

$ws.Cells.Item(33, 1).Value = 'Stephon Gilmore'
$ws.Cells.Item(33, 2).Value = 'Group2'
$ws.Cells.Item(33, 3).Value = 11
$ws.Cells.Item(33, 4).Value = 63.33333333333334
$ws.Cells.Item(33, 5).Value = 49
$ws.Cells.Item(33, 6).Value = 14.33333333333333

$ws.Cells.Item(34, 1).Value = 'Stephon Gilmore'
$ws.Cells.Item(34, 2).Value = 'Difference'
$ws.Cells.Item(34, 3).Value = 2.666666666666666
$ws.Cells.Item(34, 4).Value = 28
$ws.Cells.Item(34, 5).Value = 19.33333333333333
$ws.Cells.Item(34, 6).Value = 8.666666666666668
